# Estadisticos Matutinos 15 Oct
# Updates the mid-period statistics on "Estadisticos 1P" and "Estadisticos Final"
# (Reprobados/Aprobados/Por_Apro/Promedio), the partial "Blancos-style" count on
# "Estadisticos 2P" (column E), and refreshes the "Rescatables" (failing-students)
# roster with the current list of at-risk students.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Estadisticos 1P
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D2").Value = 12
$ws1.Range("F2").Value = 20
$ws1.Range("G2").Value = 62.5
$ws1.Range("H2").Value = 10

$ws1.Range("D3").Value = 7
$ws1.Range("F3").Value = 27
$ws1.Range("G3").Value = 79.41
$ws1.Range("H3").Value = 10

$ws1.Range("D4").Value = 10
$ws1.Range("F4").Value = 25
$ws1.Range("G4").Value = 71.43000000000001
$ws1.Range("H4").Value = 10

$ws1.Range("D5").Value = 6
$ws1.Range("F5").Value = 29
$ws1.Range("G5").Value = 82.86
$ws1.Range("H5").Value = 10

# ---------------------------------------------------------------------------
# Estadisticos 2P (only the partial "Reprobados" count is known so far)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("E2").Value = 20
$ws2.Range("E3").Value = 27
$ws2.Range("E4").Value = 25
$ws2.Range("E5").Value = 29

# ---------------------------------------------------------------------------
# Estadisticos Final (mirrors Estadisticos 1P)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("D2").Value = 12
$ws3.Range("F2").Value = 20
$ws3.Range("G2").Value = 62.5
$ws3.Range("H2").Value = 10

$ws3.Range("D3").Value = 7
$ws3.Range("F3").Value = 27
$ws3.Range("G3").Value = 79.41
$ws3.Range("H3").Value = 10

$ws3.Range("D4").Value = 10
$ws3.Range("F4").Value = 25
$ws3.Range("G4").Value = 71.43000000000001
$ws3.Range("H4").Value = 10

$ws3.Range("D5").Value = 6
$ws3.Range("F5").Value = 29
$ws3.Range("G5").Value = 82.86
$ws3.Range("H5").Value = 10

# ---------------------------------------------------------------------------
# Rescatables - refresh the roster of at-risk students (rows 2-12), and drop
# the now-obsolete trailing rows (13-17 in the old roster).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$rescatables = @(
    @(20330051920136, "MARCIAL",    "MORALES",  "IVAN DE JESUS",     "CONTRIBUYE A LA INTEGRACIÓN Y DESARROLLO DEL PERSONAL EN LA ORGANIZACIÓN", "3ARHV", 6),
    @(20330051920151, "TORRES",     "PEREZ",    "CONSTANZA XIMENA",  "CONTRIBUYE A LA INTEGRACIÓN Y DESARROLLO DEL PERSONAL EN LA ORGANIZACIÓN", "3ARHV", 6),
    @(20330051920042, "ANASTACIO",  "ROMERO",   "HIRAM FABIAN",      "ÉTICA", "3BEM", 6),
    @(19330051920162, "HUERTA",     "OFICIAL",  "MIGUEL ANGEL",      "CIENCIA, TECNOLOGÍA, SOCIEDAD Y VALORES", "5ALCM", 6),
    @(20330051920123, "GARCIA",     "JUAREZ",   "EMELIN JIROMI",     "CONTRIBUYE A LA INTEGRACIÓN Y DESARROLLO DEL PERSONAL EN LA ORGANIZACIÓN", "3ARHV", 6),
    @(20330051920129, "JIMENEZ",    "APARICIO", "YAZMIN",            "CONTRIBUYE A LA INTEGRACIÓN Y DESARROLLO DEL PERSONAL EN LA ORGANIZACIÓN", "3ARHV", 6),
    @(20330051920153, "TRUJILLO",   "",         "KIMBERLY",          "CONTRIBUYE A LA INTEGRACIÓN Y DESARROLLO DEL PERSONAL EN LA ORGANIZACIÓN", "3ARHV", 6),
    @(20330051920113, "XOTLANIHUA", "TEXCAHUA", "ALEXANDER",         "ÉTICA", "3BEM", 6),
    @(20330051920283, "ALTAMIRANO", "JUAREZ",   "KAREN ESTEPHANY",   "ÉTICA", "3BLCM", 6),
    @(20330051920390, "CARAZA",     "CRUZ",     "JARED URIEL",       "ÉTICA", "3BLCM", 6),
    @(20330051920316, "XOCUA",      "CAMPOS",   "LAURA IVETTE",      "ÉTICA", "3BLCM", 6)
)

# Remove the 5 rows (old rows 13-17) that no longer exist in the refreshed roster.
$ws4.Rows.Item(13).Resize(5).Delete()

$r = 2
foreach ($row in $rescatables) {
    $ws4.Cells.Item($r, 1).Value = $row[0]
    $ws4.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -eq "") {
        $ws4.Cells.Item($r, 3).ClearContents()
    } else {
        $ws4.Cells.Item($r, 3).Value = $row[2]
    }
    $ws4.Cells.Item($r, 4).Value = $row[3]
    $ws4.Cells.Item($r, 5).Value = $row[4]
    $ws4.Cells.Item($r, 6).Value = $row[5]
    $ws4.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
